$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.791.10'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.644.52'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').Value = '216.57'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = '0.499'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = '0.0628'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.251'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('D10').Value = '19.17'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('D11').Value = '0.0841'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').Value = '1.870.56'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = '1.651.71'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '4.16'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').Value = '64.38'
$ws.Range('E16').Value = '  -2.86%  '
$ws.Range('D17').Value = '26.812.32'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('D19').Value = '213.58'
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').Value = '4.36'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('E22').Value = '  +11.98%  '
$ws.Range('D23').Value = '6.28'
$ws.Range('E23').Value = '  -0.69%  '
$ws.Range('D24').Value = '9.36'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').Value = '145.02'
$ws.Range('E25').Value = '  -1.53%  '
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').Value = '0.118'
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').Value = '7.08'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').Value = '15.66'
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('D30').Value = '0.0510'
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('D32').Value = '3.31'
$ws.Range('E32').Value = '  -2.80%  '
$ws.Range('E33').Value = '  -1.90%  '
$ws.Range('D34').Value = '1.286.31'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').Value = '0.0173'
$ws.Range('E37').Value = '  -5.68%  '
$ws.Range('D38').Value = '0.538'
$ws.Range('D39').Value = '0.825'
$ws.Range('E39').Value = '  -0.28%  '
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').Value = '5.35'
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('D44').Value = '1.797.32'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '60.74'
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '91.36'
$ws.Range('E46').Value = '  -2.55%  '
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('E48').Value = '  -1.83%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').Value = '7.68'
$ws.Range('E50').Value = '  -2.20%  '
$ws.Range('D51').Value = '0.0977'
$ws.Range('E51').Value = '  -0.07%  '
